$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 51
$ws.Range("A51").Value = "Globo"
$ws.Range("B51").Value = "RJ TV 2"
$ws.Range("C51").Value = "Defesa Civil"
$ws.Range("D51").Value = "2025-04-04T19:13"
$ws.Range("E51").Value = "Neutro"
$ws.Range("F51").Value = "Alerta de chuvas no Norte do Estado. Em Campos e São Joao da Barra, chuvas da tarde de hoje causaram transtornos. Repórter *ao vivo*. Não chove no momento. Temperaturas caíram. Hoje à tarde, Prefeito Wladimir Garotinho divulgou vídeo nas redes sociais com representante do Centro de Monitoramento de Desastres da Defesa Civil tranquilizando a população e informando que a situação é melhor hoje. "

# Row 52
$ws.Range("A52").Value = "Globo"
$ws.Range("B52").Value = "RJ TV 2"
$ws.Range("C52").Value = "Defesa Civil"
$ws.Range("D52").Value = "2025-04-04T19:22"
$ws.Range("E52").Value = "Positivo"
$ws.Range("F52").Value = "Alerta de chuvas no Norte. Defesas Civis de Campos, São João da Barra e São Fidélis estão em alerta. Equipes da Defesa Civil de Campos estão em estágio de atenção. Entrevista com o subsecretário da Defesa Civil, major Edison Pessanha. Previsão de chuva para o final de semana. Outros municípios também foram citados. *matéria*"

$wb.Save()
